# Updated symbol list - apply cell-level edits from the crypto price refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.82"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.22%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.87%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.225"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.10%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05917"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.88%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.687"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.15%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8700"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.30%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.072"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'25.02%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1417"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.52%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.26%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.03238"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.25%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09237"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.28%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001548"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.42%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0006065"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.49%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005993"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.64%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.484"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.12%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.270"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Value = "'0.3151"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.57%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.03626"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'11.15%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1293"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.47%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.557"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.89%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04185"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.39%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.1403"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.71%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.05%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004536"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'9.55%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.16%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001943"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'34.18%"
$ws.Range("E28").Style = "Normal"
$ws.Range("E40").Value = "'1.35%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1106"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.96%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").Value = "'0.003993"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-22.89%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002386"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-0.62%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009957"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'8.80%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005449"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'3.10%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000752"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.30%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.1094"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'4.23%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-12.11%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.30%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.30%"
$ws.Range("E50").Style = "Normal"
